# edit.ps1
# Applies the "Update countries & provincias Spain" edit to paises.xlsx:
#  - Refreshes COVID-19 case figures for several countries (new data pulled
#    at 20:30 instead of 19:13), which re-sorts the table (descending by
#    "Casos totales") and causes a handful of adjacent-rank countries to
#    swap row positions.
#  - Updates the "Datos actualizados ..." timestamp cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 6496457
$ws.Range("C4").Value = 10882
$ws.Range("D4").Value = 3768762
$ws.Range("E4").Value = 2533974
$ws.Range("G4").Value = 187
$ws.Range("H4").Value = 193721
# Row 5
$ws.Range("B5").Value = 4363843
$ws.Range("C5").Value = 86259
$ws.Range("D5").Value = 3393583
$ws.Range("E5").Value = 896366
$ws.Range("G5").Value = 1078
$ws.Range("H5").Value = 73894
# Row 12
$ws.Range("B12").Value = 534513
$ws.Range("C12").Value = 8964
$ws.Range("G12").Value = 78
$ws.Range("H12").Value = 29594
# Row 16
$ws.Range("B16").Value = 352560
$ws.Range("C16").Value = 2460
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = 41586
# Row 17
$ws.Range("A17").Value = "Francia"
$ws.Range("B17").Value = 335524
$ws.Range("C17").Value = 6544
$ws.Range("D17").Value = 88226
$ws.Range("E17").Value = 216572
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 30726
# Row 18
$ws.Range("A18").Value = "Banglades"
$ws.Range("B18").Value = 329251
$ws.Range("C18").Value = 1892
$ws.Range("D18").Value = 227809
$ws.Range("E18").Value = 96890
$ws.Range("G18").Value = 36
$ws.Range("H18").Value = 4552
# Row 45
$ws.Range("A45").Value = "Marruecos"
$ws.Range("B45").Value = 75721
$ws.Range("C45").Value = 1941
$ws.Range("D45").Value = 57239
$ws.Range("E45").Value = 17055
$ws.Range("G45").Value = 33
$ws.Range("H45").Value = 1427
# Row 46
$ws.Range("A46").Value = "Emiratos Arabes Unidos"
$ws.Range("B46").Value = 75098
$ws.Range("C46").Value = 644
$ws.Range("D46").Value = 66943
$ws.Range("E46").Value = 7764
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 391
# Row 52
$ws.Range("B52").Value = 60784
$ws.Range("C52").Value = 1136
$ws.Range("D52").Value = 22677
$ws.Range("E52").Value = 37158
$ws.Range("G52").Value = 16
$ws.Range("H52").Value = 949
# Row 81
$ws.Range("B81").Value = 19583
$ws.Range("C81").Value = 749
$ws.Range("D81").Value = 2247
$ws.Range("E81").Value = 17022
$ws.Range("G81").Value = 18
$ws.Range("H81").Value = 314
# Row 93
$ws.Range("E93").Value = 3993
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 321
# Row 99
$ws.Range("A99").Value = "Namibia"
$ws.Range("B99").Value = 8928
$ws.Range("C99").Value = 118
$ws.Range("D99").Value = 3981
$ws.Range("E99").Value = 4856
$ws.Range("H99").Value = 91
# Row 100
$ws.Range("A100").Value = "Tayikistan"
$ws.Range("B100").Value = 8860
$ws.Range("C100").Value = 36
$ws.Range("D100").Value = 7650
$ws.Range("E100").Value = 1140
$ws.Range("H100").Value = 70
# Row 132
$ws.Range("A132").Value = "Sri Lanka"
$ws.Range("B132").Value = 3140
$ws.Range("C132").Value = 17
$ws.Range("D132").Value = 2935
$ws.Range("E132").Value = 193
$ws.Range("H132").Value = 12
# Row 133
$ws.Range("A133").Value = "Lituania"
$ws.Range("C133").Value = 31
$ws.Range("D133").Value = 1994
$ws.Range("E133").Value = 1051
$ws.Range("H133").Value = 86
# Row 135
$ws.Range("B135").Value = 2882
$ws.Range("C135").Value = 12
$ws.Range("D135").Value = 2258
# Row 137
$ws.Range("A137").Value = "Jordania"
$ws.Range("B137").Value = 2581
$ws.Range("C137").Value = 103
$ws.Range("D137").Value = 1885
$ws.Range("E137").Value = 677
$ws.Range("G137").Value = 2
$ws.Range("H137").Value = 19
# Row 138
$ws.Range("A138").Value = "Estonia"
$ws.Range("B138").Value = 2564
$ws.Range("C138").Value = 32
$ws.Range("D138").Value = 2195
$ws.Range("E138").Value = 305
$ws.Range("H138").Value = 64
# Row 139
$ws.Range("A139").Value = "Sudan del Sur"
$ws.Range("B139").Value = 2545
$ws.Range("D139").Value = 1290
$ws.Range("E139").Value = 1207
$ws.Range("H139").Value = 48
# Row 140
$ws.Range("A140").Value = "Aruba"
$ws.Range("B140").Value = 2482
$ws.Range("D140").Value = 1244
$ws.Range("E140").Value = 1223
$ws.Range("H140").Value = 15
# Row 145
$ws.Range("B145").Value = 2150
$ws.Range("C145").Value = 7
$ws.Range("D145").Value = 2060
$ws.Range("E145").Value = 80
# Row 149
$ws.Range("B149").Value = 1994
$ws.Range("C149").Value = 5
$ws.Range("D149").Value = 1203
$ws.Range("G149").Value = 3
$ws.Range("H149").Value = 576

# Timestamp header cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 20:30"
